$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 40
$ws.Range("H40").Value = 983.75
$ws.Range("I40").Value = 920
$ws.Range("J40").Value = 1033.3334
$ws.Range("K40").Value = 920
$ws.Range("L40").Value = 1033.3334
$ws.Range("M40").Value = -745
$ws.Range("N40").Value = -1383.3334

# ALC row 113
$ws.Range("H113").Value = 1852.6154
$ws.Range("I113").Value = 1851
$ws.Range("J113").Value = 1853.3334
$ws.Range("K113").Value = 1851
$ws.Range("L113").Value = 1853.3334
$ws.Range("M113").Value = 1403
$ws.Range("N113").Value = -8361.3334

# ALC row 123
$ws.Range("H123").Value = 33593.35
$ws.Range("I123").Value = 26983.857
$ws.Range("J123").Value = 37152.31
$ws.Range("K123").Value = 26983.857
$ws.Range("L123").Value = 37152.31
$ws.Range("M123").Value = -22083.857
$ws.Range("N123").Value = -46952.31

# ALC row 129
$ws.Range("H129").Value = 799.7778
$ws.Range("I129").Value = 354.85715
$ws.Range("J129").Value = 907.1724
$ws.Range("K129").Value = 1064.57145
$ws.Range("L129").Value = 2721.5172
$ws.Range("M129").Value = 3935.42855
$ws.Range("N129").Value = -12721.5172

# ALC row 132
$ws.Range("H132").Value = 3714.8386
$ws.Range("I132").Value = 3805.5173
$ws.Range("J132").Value = 2400
$ws.Range("K132").Value = 11416.5519
$ws.Range("L132").Value = 7200
$ws.Range("M132").Value = -8886.5519
$ws.Range("N132").Value = -12260

# ALC row 141
$ws.Range("H141").Value = 3027.7273
$ws.Range("I141").Value = 1306.875
$ws.Range("J141").Value = 7616.6665
$ws.Range("K141").Value = 3920.625
$ws.Range("L141").Value = 22849.9995
$ws.Range("M141").Value = 1259.375
$ws.Range("N141").Value = -33209.99950000001

$ws = $wb.Worksheets.Item("ARM")
# ARM row 2
$ws.Range("H2").Value = 1963156.9
$ws.Range("I2").Value = 1832
$ws.Range("J2").Value = 2943819.2
$ws.Range("K2").Value = 1832
$ws.Range("L2").Value = 2943819.2
$ws.Range("M2").Value = -1719
$ws.Range("N2").Value = -2944045.2

# ARM row 32
$ws.Range("H32").Value = 10352.426
$ws.Range("I32").Value = 8285.731
$ws.Range("J32").Value = 24474.834
$ws.Range("K32").Value = 8285.731
$ws.Range("L32").Value = 24474.834
$ws.Range("M32").Value = -7998.731
$ws.Range("N32").Value = -25048.834

# ARM row 116
$ws.Range("H116").Value = 1963156.9
$ws.Range("I116").Value = 1832
$ws.Range("J116").Value = 2943819.2
$ws.Range("K116").Value = 1832
$ws.Range("L116").Value = 2943819.2
$ws.Range("M116").Value = 462
$ws.Range("N116").Value = -2948407.2

# ARM row 132
$ws.Range("H132").Value = 24944.432
$ws.Range("I132").Value = 2205.8647
$ws.Range("J132").Value = 145134
$ws.Range("K132").Value = 6617.5941
$ws.Range("L132").Value = 435402
$ws.Range("M132").Value = -4087.5941
$ws.Range("N132").Value = -440462

$ws = $wb.Worksheets.Item("BSM")
# BSM row 3
$ws.Range("H3").Value = 1963156.9
$ws.Range("I3").Value = 1832
$ws.Range("J3").Value = 2943819.2
$ws.Range("K3").Value = 1832
$ws.Range("L3").Value = 2943819.2
$ws.Range("M3").Value = -1718
$ws.Range("N3").Value = -2944047.2

$ws = $wb.Worksheets.Item("CRP")
# CRP row 16
$ws.Range("H16").Value = 1553.2916
$ws.Range("I16").Value = 1756.7894
$ws.Range("J16").Value = 780
$ws.Range("K16").Value = 1756.7894
$ws.Range("L16").Value = 780
$ws.Range("M16").Value = -1469.7894
$ws.Range("N16").Value = -1354

# CRP row 59
$ws.Range("H59").Value = 14381.75
$ws.Range("I59").Value = 10500
$ws.Range("J59").Value = 15158.1
$ws.Range("K59").Value = 10500
$ws.Range("L59").Value = 15158.1
$ws.Range("M59").Value = -9355
$ws.Range("N59").Value = -17448.1

# CRP row 68
$ws.Range("H68").Value = 18733.666
$ws.Range("J68").Value = 18733.666
$ws.Range("L68").Value = 18733.666
$ws.Range("N68").Value = -20231.666

# CRP row 71
$ws.Range("H71").Value = 18733.666
$ws.Range("J71").Value = 18733.666
$ws.Range("L71").Value = 56200.99800000001
$ws.Range("N71").Value = -63688.99800000001

# CRP row 107
$ws.Range("H107").Value = 2000.6
$ws.Range("I107").Value = 2408.0715
$ws.Range("J107").Value = 1049.8334
$ws.Range("K107").Value = 2408.0715
$ws.Range("L107").Value = 1049.8334
$ws.Range("M107").Value = -488.0715
$ws.Range("N107").Value = -4889.8334

# CRP row 113
$ws.Range("H113").Value = 1553.2916
$ws.Range("I113").Value = 1756.7894
$ws.Range("J113").Value = 780
$ws.Range("K113").Value = 1756.7894
$ws.Range("L113").Value = 780
$ws.Range("M113").Value = 413.2106000000001
$ws.Range("N113").Value = -5120

# CRP row 125
$ws.Range("H125").Value = 50653.332
$ws.Range("J125").Value = 50653.332
$ws.Range("L125").Value = 50653.332
$ws.Range("N125").Value = -55573.332

$ws = $wb.Worksheets.Item("CUL")
# CUL row 122
$ws.Range("H122").Value = 12346522
$ws.Range("I122").Value = 15152013
$ws.Range("K122").Value = 136368117
$ws.Range("M122").Value = -136365667

# CUL row 129
$ws.Range("H129").Value = 1333.7037
$ws.Range("I129").Value = 796.55554
$ws.Range("J129").Value = 1602.2778
$ws.Range("K129").Value = 2389.66662
$ws.Range("L129").Value = 4806.8334
$ws.Range("M129").Value = 2610.33338
$ws.Range("N129").Value = -14806.8334

# CUL row 137
$ws.Range("H137").Value = 24284.844
$ws.Range("J137").Value = 29385.537
$ws.Range("L137").Value = 88156.611
$ws.Range("N137").Value = -98356.611

$ws = $wb.Worksheets.Item("GSM")
# GSM row 107
$ws.Range("H107").Value = 388.3871
$ws.Range("I107").Value = 348
$ws.Range("J107").Value = 461.81818
$ws.Range("K107").Value = 348
$ws.Range("L107").Value = 461.81818
$ws.Range("M107").Value = 1572
$ws.Range("N107").Value = -4301.81818

# GSM row 113
$ws.Range("H113").Value = 1133.6
$ws.Range("I113").Value = 1162.2858
$ws.Range("J113").Value = 1066.6666
$ws.Range("K113").Value = 1162.2858
$ws.Range("L113").Value = 1066.6666
$ws.Range("M113").Value = 1007.7142
$ws.Range("N113").Value = -5406.6666

# GSM row 124
$ws.Range("H124").Value = 54840
$ws.Range("J124").Value = 54840
$ws.Range("L124").Value = 54840
$ws.Range("N124").Value = -64660

# GSM row 126
$ws.Range("H126").Value = 1921.625
$ws.Range("I126").Value = 1493.75
$ws.Range("J126").Value = 2349.5
$ws.Range("K126").Value = 4481.25
$ws.Range("L126").Value = 7048.5
$ws.Range("M126").Value = -2011.25
$ws.Range("N126").Value = -11988.5

$ws = $wb.Worksheets.Item("WVR")
# WVR row 126
$ws.Range("H126").Value = 2231.5
$ws.Range("I126").Value = 1891.2273
$ws.Range("J126").Value = 3479.1667
$ws.Range("K126").Value = 5673.6819
$ws.Range("L126").Value = 10437.5001
$ws.Range("M126").Value = -3203.6819
$ws.Range("N126").Value = -15377.5001

# WVR row 132
$ws.Range("H132").Value = 4502.979
$ws.Range("I132").Value = 6583.357
$ws.Range("J132").Value = 1590.45
$ws.Range("K132").Value = 19750.071
$ws.Range("L132").Value = 4771.35
$ws.Range("M132").Value = -17220.071
$ws.Range("N132").Value = -9831.35
